$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I ("I0") and J ("IF"), styled like the other headers (copy format from H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..55 (columns I and J)
$I = @(6, 6, 5, 8, 6, 7, 9, 6, 5, 8, 11, 6, 6, 10, 7, 5, 7, 6, 8, 7, 7, 6, 7, 8, 6, 6, 9, 7, 8, 7, 7, 6, 9, 9, 8, 9, 8, 8, 9, 8, 8, 8, 7, 4, 9, 8, 8, 7, 8, 4, 4, 4, 4, 6)
$J = @(6, 6, 5, 8, 6, 8, 9, 6, 5, 8, 11, 6, 7, 10, 7, 5, 7, 6, 8, 7, 8, 6, 7, 8, 6, 6, 9, 7, 8, 8, 7, 6, 9, 9, 8, 9, 8, 8, 9, 8, 8, 8, 7, 4, 9, 8, 8, 7, 8, 4, 4, 4, 4, 6)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}
